# Update the "F" column (numeric count) values on the "展览" and
# "全部类型" worksheets. Both sheets contain the same underlying rows,
# but the totals were recomputed independently, so the "全部类型" sheet
# gets a slightly different value for row 40 (3380) than "展览" (3379).

$wb = $excel.ActiveWorkbook

# row -> new F value, shared by both sheets
$commonUpdates = @{
    2  = 39
    7  = 106
    8  = 74
    9  = 442
    12 = 546
    14 = 282
    15 = 21
    16 = 341
    18 = 86
    20 = 42
    21 = 86
    22 = 843
    23 = 1366
    24 = 287
    25 = 299
    27 = 62
    29 = 33
    30 = 81
    31 = 200
    33 = 256
    34 = 1587
    35 = 46
    38 = 560
    41 = 396
    43 = 868
    45 = 56
    46 = 39
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $commonUpdates.Keys) {
    $ws1.Range("F$row").Value = $commonUpdates[$row]
}
$ws1.Range("F40").Value = 3379

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $commonUpdates.Keys) {
    $ws4.Range("F$row").Value = $commonUpdates[$row]
}
$ws4.Range("F40").Value = 3380
